# Refresh "t11.2" sheet: data now covers 2003-2023 (one fewer year than
# before - the old 2002 row is dropped, every row shifts up one year, and
# the final row gets freshly recomputed 2023 figures). Two cells (H5/G10)
# now hold the literal "-" instead of a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2003
$ws.Range("B2").Value = -22.385839805274198
$ws.Range("C2").Value = 9.2603186402582018
$ws.Range("D2").Value = 14.954083940473083
$ws.Range("E2").Value = 2.3873272486277575
$ws.Range("F2").Value = -37.662446672212425
$ws.Range("G2").Value = 131.80246779981033
$ws.Range("H2").Value = -7.9143785445388026
$ws.Range("I2").Value = -8.8749697009830051

$ws.Range("A3").Value = 2004
$ws.Range("B3").Value = 21.114475382626964
$ws.Range("C3").Value = 9.2439699782845164
$ws.Range("D3").Value = 19.822242732639616
$ws.Range("E3").Value = -71.815596323481017
$ws.Range("F3").Value = 409.76452479925882
$ws.Range("G3").Value = 1149.8301898139089
$ws.Range("H3").Value = 327.32329405592128
$ws.Range("I3").Value = 46.654274071015102

$ws.Range("A4").Value = 2005
$ws.Range("B4").Value = 43.182225331122922
$ws.Range("C4").Value = 12.422991544721818
$ws.Range("D4").Value = 145.51136896807009
$ws.Range("E4").Value = -68.356303257987975
$ws.Range("F4").Value = 28.173654920046975
$ws.Range("G4").Value = 29.387731767407054
$ws.Range("H4").Value = -99.917476977043478
$ws.Range("I4").Value = 23.745148659270001

$ws.Range("A5").Value = 2006
$ws.Range("B5").Value = 1.2287404158958015
$ws.Range("C5").Value = 21.83023966606601
$ws.Range("D5").Value = -29.134586891372017
$ws.Range("E5").Value = 1232.1819603235679
$ws.Range("F5").Value = -74.12778066924966
$ws.Range("G5").Value = 94.85649407070072
$ws.Range("H5").Value = "-"
$ws.Range("I5").Value = -16.695188310612984

$ws.Range("A6").Value = 2007
$ws.Range("B6").Value = 8.7442263739687895
$ws.Range("C6").Value = -5.3900352548983843
$ws.Range("D6").Value = -19.025904057174202
$ws.Range("E6").Value = 21.963731344983749
$ws.Range("F6").Value = 0.76629024777563615
$ws.Range("G6").Value = -97.098284268128864
$ws.Range("H6").Value = -8.7824427786134756
$ws.Range("I6").Value = 3.6660157288067152

$ws.Range("A7").Value = 2008
$ws.Range("B7").Value = 9.7650286850176737
$ws.Range("C7").Value = 4.7904909817143126
$ws.Range("D7").Value = 12.53670027258087
$ws.Range("E7").Value = 6.9380172169372978
$ws.Range("F7").Value = -2.1202467170273298
$ws.Range("G7").Value = 74.496830320100145
$ws.Range("H7").Value = 13.169589571007535
$ws.Range("I7").Value = 5.6507261321199875

$ws.Range("A8").Value = 2009
$ws.Range("B8").Value = 11.426940125995699
$ws.Range("C8").Value = 5.833801714606146
$ws.Range("D8").Value = 29.189303582080917
$ws.Range("E8").Value = -3.0431472113520486
$ws.Range("F8").Value = -22.716336551199646
$ws.Range("G8").Value = -50.393422220465879
$ws.Range("H8").Value = 26.114779710459256
$ws.Range("I8").Value = 1.8754414460817381

$ws.Range("A9").Value = 2010
$ws.Range("B9").Value = 17.093975538967921
$ws.Range("C9").Value = 40.238473949626027
$ws.Range("D9").Value = 1.0429042330415639
$ws.Range("E9").Value = 2.5458731551389224
$ws.Range("F9").Value = 2.5231548646748037
$ws.Range("G9").Value = -99.974945001688411
$ws.Range("H9").Value = -0.57123164564236717
$ws.Range("I9").Value = 20.238216602530869

$ws.Range("A10").Value = 2011
$ws.Range("B10").Value = 5.6076524400949568
$ws.Range("C10").Value = 6.0935625969827267
$ws.Range("D10").Value = -0.83415743431578537
$ws.Range("E10").Value = 2.1979525682162437
$ws.Range("F10").Value = 70.000994707275225
$ws.Range("G10").Value = "-"
$ws.Range("H10").Value = 7.8590308793093033
$ws.Range("I10").Value = 10.733017285221514

$ws.Range("A11").Value = 2012
$ws.Range("B11").Value = 0.037933974584225538
$ws.Range("C11").Value = 7.1777907518893347
$ws.Range("D11").Value = -1.8934740050543963
$ws.Range("E11").Value = 5.8182612285416546
$ws.Range("F11").Value = -32.523727770598065
$ws.Range("G11").Value = 23.018467490521388
$ws.Range("H11").Value = -17.809975639150011
$ws.Range("I11").Value = -0.55679798168963845

$ws.Range("A12").Value = 2013
$ws.Range("B12").Value = 5.5017854241490083
$ws.Range("C12").Value = -10.653822785987133
$ws.Range("D12").Value = -5.052178173926702
$ws.Range("E12").Value = 1.1001059591388307
$ws.Range("F12").Value = 3.0269421172712363
$ws.Range("G12").Value = 7.7796436384257373
$ws.Range("H12").Value = -20.553987932815289
$ws.Range("I12").Value = -2.9398019518870888

$ws.Range("A13").Value = 2014
$ws.Range("B13").Value = 1.4172464954842479
$ws.Range("C13").Value = 4.591953361346568
$ws.Range("D13").Value = -0.69087327871165671
$ws.Range("E13").Value = -8.7302493009758226
$ws.Range("F13").Value = -20.420474315066571
$ws.Range("G13").Value = -8.7260252866161565
$ws.Range("H13").Value = -15.609979096097037
$ws.Range("I13").Value = -1.4676473703888915

$ws.Range("A14").Value = 2015
$ws.Range("B14").Value = -6.7066481104024245
$ws.Range("C14").Value = -18.547082361468803
$ws.Range("D14").Value = -6.744644356555507
$ws.Range("E14").Value = -9.4442500277878594
$ws.Range("F14").Value = 41.401598569669098
$ws.Range("G14").Value = 10.403799081898569
$ws.Range("H14").Value = 22.262826999968176
$ws.Range("I14").Value = -8.4534001435952142

$ws.Range("A15").Value = 2016
$ws.Range("B15").Value = 5.0868525323150227
$ws.Range("C15").Value = -26.316075234334079
$ws.Range("D15").Value = -0.86357487054186111
$ws.Range("E15").Value = 1.5933051991474656
$ws.Range("F15").Value = 51.007190235015543
$ws.Range("G15").Value = -1.7256676403272553
$ws.Range("H15").Value = 14.708363135761848
$ws.Range("I15").Value = -2.0534001149640435

$ws.Range("A16").Value = 2017
$ws.Range("B16").Value = 4.836710328016558
$ws.Range("C16").Value = 13.877031874092305
$ws.Range("D16").Value = 5.074170390703614
$ws.Range("E16").Value = -4.1938469448560278
$ws.Range("F16").Value = -0.1927224985693865
$ws.Range("G16").Value = 5.5633593564111905
$ws.Range("H16").Value = -20.128763147442687
$ws.Range("I16").Value = 4.8342380071453972

$ws.Range("A17").Value = 2018
$ws.Range("B17").Value = 4.024003573824908
$ws.Range("C17").Value = -4.0221664623212199
$ws.Range("D17").Value = 18.723798651272428
$ws.Range("E17").Value = -17.365704431327657
$ws.Range("F17").Value = -1.5906640669415606
$ws.Range("G17").Value = 2.6681716367631614
$ws.Range("H17").Value = 39.004405147659881
$ws.Range("I17").Value = -2.7001156815860683

$ws.Range("A18").Value = 2019
$ws.Range("B18").Value = 1.4693904127345947
$ws.Range("C18").Value = 3.4104839329791581
$ws.Range("D18").Value = 9.262331720820228
$ws.Range("E18").Value = -11.894857852363627
$ws.Range("F18").Value = -11.614068249684006
$ws.Range("G18").Value = -12.058649059527193
$ws.Range("H18").Value = 3.2210663880971602
$ws.Range("I18").Value = -1.9726738289750689

$ws.Range("A19").Value = 2020
$ws.Range("B19").Value = 12.114229720561575
$ws.Range("C19").Value = -11.382717198232539
$ws.Range("D19").Value = 3.5780148003290213
$ws.Range("E19").Value = -3.8473115346673747
$ws.Range("F19").Value = 5.8026948396254996
$ws.Range("G19").Value = 12.121071464460442
$ws.Range("H19").Value = 19.113815485362263
$ws.Range("I19").Value = 1.0652294522354167

$ws.Range("A20").Value = 2021
$ws.Range("B20").Value = 14.930988249598576
$ws.Range("C20").Value = 11.572202695009469
$ws.Range("D20").Value = 14.605859836417622
$ws.Range("E20").Value = -6.7563351590502529
$ws.Range("F20").Value = 12.992686582931201
$ws.Range("G20").Value = -13.515223422290445
$ws.Range("H20").Value = 40.011753552481764
$ws.Range("I20").Value = 10.983691226005021

$ws.Range("A21").Value = 2022
$ws.Range("B21").Value = 8.5320655032458568
$ws.Range("C21").Value = 0.60545028858061123
$ws.Range("D21").Value = 7.6822467147618978
$ws.Range("E21").Value = -20.425447152621157
$ws.Range("F21").Value = 28.880669369150034
$ws.Range("G21").Value = 47.239280579476485
$ws.Range("H21").Value = -25.594183192041719
$ws.Range("I21").Value = 6.8397493130987863

$ws.Range("A22").Value = 2023
$ws.Range("B22").Value = -51.383979806680514
$ws.Range("C22").Value = -56.80633458157287
$ws.Range("D22").Value = -52.815975147785551
$ws.Range("E22").Value = -67.782323122617768
$ws.Range("F22").Value = -55.910515333732633
$ws.Range("G22").Value = -58.596883880607486
$ws.Range("H22").Value = -44.206694596859677
$ws.Range("I22").Value = -55.115612356644895

# Drop the now-surplus last row (old 2023 data, already superseded above).
$ws.Rows("23:23").Delete()

# Match the saved selection state (whole used range selected).
$ws.Range("A1:I22").Select()
